# Rename the inline picture objects (the wp:docPr / pic:cNvPr "name"
# attribute) that are embedded in the document's headers and footers.
#
#   header (primary + first page) : BTec logo   image2.jpg -> image1.jpg
#   footer (primary + first page) : Pearson logo image1.png -> image2.png
#
# InlineShape has no writable .Name property in the Word object model
# (only the floating Shape object does), so each picture is briefly
# converted to a Shape, renamed, and converted back to an InlineShape -
# this mirrors how Word itself would have to perform the rename and
# leaves the picture inline (wp:inline) exactly as before.

function Rename-InlinePicture($range, $newName) {
    $inlineShape = $range.InlineShapes.Item(1)
    $shape = $inlineShape.ConvertToShape()
    $shape.Name = $newName
    $shape.ConvertToInlineShape() | Out-Null
}

$d = $word.ActiveDocument
$section = $d.Sections.Item(1)

# Headers: BTec_Logo-Orange picture, image2.jpg -> image1.jpg
Rename-InlinePicture $section.Headers.Item(1).Range "image1.jpg"
Rename-InlinePicture $section.Headers.Item(2).Range "image1.jpg"

# Footers: Pearson logo picture, image1.png -> image2.png
Rename-InlinePicture $section.Footers.Item(1).Range "image2.png"
Rename-InlinePicture $section.Footers.Item(2).Range "image2.png"
